# ControlsTable.xlsx - "End of 114 video. Add basic tank camera controls."
#
# Sheet2 currently ends at row 11 with the old Base-Rotate / Turret-Rotate
# rows (6-11). We rework rows 6-11 (Base rotate -> Body rotate / Turret
# controls) and append new rows 12-15 for Barrel elevation controls.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# xlCenter = -4108, xlBottom = -4107

# --- Row 6/7 block: "Rotate left" -> "Rotate Body Left" -------------------
$ws.Range("A6").Value = "Rotate Body Left"

$ws.Range("B6").Value = "A-Sym Triggers/`nBumpers"
$ws.Range("B6").WrapText = $true

$ws.Range("C6").Value = "Left Stick Left"
$ws.Range("C7").Value = "A Key"

$ws.Range("D6").Value = "Tracks - Different Speeds"

# Row 6 grew a wrapped, 2-line value -- Excel auto-fit bumps the row
# height; AutoFit the row back down so no stray ht/customHeight sticks.
$ws.Rows.Item(6).AutoFit()

# --- Row 8/9 block: "Base Rotate right" -> "Turret Rotate Left" -----------
$ws.Range("A8").Value = "Turret Rotate Left"

$rB89 = $ws.Range("B8:B9")
$rB89.ClearContents()
$rB89.HorizontalAlignment = -4108
$rB89.VerticalAlignment = -4107
$rB89.WrapText = $false

$ws.Range("C8").Value = "Right Stick Left"
$ws.Range("C9").Value = "Mouse Left"

$ws.Range("D8").Value = "Turret Rotate Left"

# --- Row 10/11 block (new): "Turret Rotate Right" --------------------------
$rA1011 = $ws.Range("A10:A11")
$rA1011.HorizontalAlignment = -4108
$rA1011.VerticalAlignment = -4108
$rA1011.WrapText = $false
$rA1011.MergeCells = $true
$ws.Range("A10").Value = "Turret Rotate Right"

$rB1011 = $ws.Range("B10:B11")
$rB1011.HorizontalAlignment = -4108
$rB1011.VerticalAlignment = -4107
$rB1011.WrapText = $false
$rB1011.MergeCells = $true

$ws.Range("C10").Value = "Right Stick Right"
$ws.Range("C10").VerticalAlignment = -4108
$ws.Range("C11").Value = "Mouse Right"
$ws.Range("C11").VerticalAlignment = -4108

$rD1011 = $ws.Range("D10:D11")
$rD1011.HorizontalAlignment = -4108
$rD1011.VerticalAlignment = -4108
$rD1011.WrapText = $false
$rD1011.MergeCells = $true
$ws.Range("D10").Value = "Turret Rotate Right"

# --- Row 12/13 block (new): "Barrel Up" ------------------------------------
$rA1213 = $ws.Range("A12:A13")
$rA1213.HorizontalAlignment = -4108
$rA1213.VerticalAlignment = -4108
$rA1213.WrapText = $false
$rA1213.MergeCells = $true
$ws.Range("A12").Value = "Barrel Up"

$rB1213 = $ws.Range("B12:B13")
$rB1213.HorizontalAlignment = -4108
$rB1213.VerticalAlignment = -4108
$rB1213.WrapText = $false
$rB1213.MergeCells = $true

$ws.Range("C12").Value = "Right Stick Up"
$ws.Range("C12").VerticalAlignment = -4108
$ws.Range("C13").Value = "Mouse Up"
$ws.Range("C13").VerticalAlignment = -4108

$rD1213 = $ws.Range("D12:D13")
$rD1213.HorizontalAlignment = -4108
$rD1213.VerticalAlignment = -4108
$rD1213.WrapText = $false
$rD1213.MergeCells = $true
$ws.Range("D12").Value = "Turret Rotator"

# --- Row 14/15 block (new): "Barrel Down" ----------------------------------
$rA1415 = $ws.Range("A14:A15")
$rA1415.HorizontalAlignment = -4108
$rA1415.VerticalAlignment = -4108
$rA1415.WrapText = $false
$rA1415.MergeCells = $true
$ws.Range("A14").Value = "Barrel Down"

$rB1415 = $ws.Range("B14:B15")
$rB1415.HorizontalAlignment = -4108
$rB1415.VerticalAlignment = -4108
$rB1415.WrapText = $false
$rB1415.MergeCells = $true

$ws.Range("C14").Value = "Right Stick Down"
$ws.Range("C14").VerticalAlignment = -4108
$ws.Range("C15").Value = "Mouse  Down"
$ws.Range("C15").VerticalAlignment = -4108

$rD1415 = $ws.Range("D14:D15")
$rD1415.HorizontalAlignment = -4108
$rD1415.VerticalAlignment = -4108
$rD1415.WrapText = $false
$rD1415.MergeCells = $true
$ws.Range("D14").Value = "Barrel Elevator"

# --- Selection / active cell housekeeping (matches saved view state) ------
$ws.Activate()
$ws.Range("D16").Select()
